$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "master": insert a new row 2 for the "Gifts for under £20 cta"
# link, pushing the existing rows (cta / Left / Right) down by one.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("master")

$ws1.Rows("2:2").Insert()

$ws1.Range("A2").Value = "Gifts for under £20 cta"
$ws1.Range("B2").Value = "cta"
$ws1.Range("B2").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Sheet "carousels": remove the "Browse gifts block 7/8" rows and
# replace the "Gifts for under £40" rows 3-7 with a new "Gifts for
# under £20" module (blocks 1-4) plus "Gifts for under £40" blocks 1-2,
# renumbering the remaining "Gifts for under £40" rows to blocks 3-4.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("carousels")

# Drop the two "Browse gifts block 7/8" rows (old rows 7 and 8).
$ws2.Rows("7:8").Delete()

# Make room for the 8 rows the new module needs (currently only 5
# remain: the old "Gifts for under £40 block 3-7" rows).
$ws2.Rows("7:9").Insert()

$ws2.Range("A7").Value = "Gifts for under £20 block 1"
$ws2.Range("B7").Value = "gfu201"
$ws2.Range("B7").Style = "Hyperlink"

$ws2.Range("A8").Value = "Gifts for under £20 block 2"
$ws2.Range("B8").Value = "gfu202"
$ws2.Range("B8").Style = "Hyperlink"

$ws2.Range("A9").Value = "Gifts for under £20 block 3"
$ws2.Range("B9").Value = "gfu203"
$ws2.Range("B9").Style = "Hyperlink"

$ws2.Range("A10").Value = "Gifts for under £20 block 4"
$ws2.Range("B10").Value = "gfu204"
$ws2.Range("B10").Style = "Hyperlink"

$ws2.Range("A11").Value = "Gifts for under £40 block 1"
$ws2.Range("B11").Value = "gfu401"
$ws2.Range("B11").Style = "Hyperlink"

$ws2.Range("A12").Value = "Gifts for under £40 block 2"
$ws2.Range("B12").Value = "gfu402"
$ws2.Range("B12").Style = "Hyperlink"

$ws2.Range("A13").Value = "Gifts for under £40 block 3"
$ws2.Range("B13").Value = "gfu403"
$ws2.Range("B13").Style = "Hyperlink"

$ws2.Range("A14").Value = "Gifts for under £40 block 4"
$ws2.Range("B14").Value = "gfu404"
$ws2.Range("B14").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Restore sensible selections on both sheets (cosmetic, matches the
# state Excel leaves the workbook in after these edits).
# ---------------------------------------------------------------------
$ws1.Range("A3").Select()
$ws2.Range("I6").Select()

$wb.Save()
